{"js": "const replacements = [\n  [\"87\u00d733=2871\", \"96\u00d721=2016\"],\n  [\"60\u00d773=4380\", \"73\u00d788=6424\"],\n  [\"51\u00d784=4284\", \"61\u00d742=2562\"],\n  [\"84\u00d727=2268\", \"56\u00d779=4424\"],\n  [\"67\u00d744=2948\", \"33\u00d712=396\"],\n  [\"36\u00d751=1836\", \"68\u00d776=5168\"],\n  [\"13\u00d737=481\", \"62\u00d745=2790\"],\n  [\"41\u00d765=2665\", \"51\u00d739=1989\"],\n  [\"81\u00d738=3078\", \"39\u00d780=3120\"],\n  [\"51\u00d717=867\", \"53\u00d713=689\"],\n  [\"31\u00d714=434\", \"72\u00d733=2376\"],\n  [\"22\u00d744=968\", \"83\u00d759=4897\"],\n  [\"43\u00d716=688\", \"49\u00d788=4312\"],\n  [\"84\u00d784=7056\", \"66\u00d731=2046\"],\n  [\"67\u00d775=5025\", \"16\u00d771=1136\"],\n  [\"50\u00d713=650\", \"21\u00d774=1554\"],\n  [\"87\u00d779=6873\", \"88\u00d764=5632\"],\n  [\"33\u00d780=2640\", \"27\u00d789=2403\"],\n  [\"91\u00d721=1911\", \"61\u00d723=1403\"],\n  [\"64\u00d711=704\", \"84\u00d734=2856\"],\n  [\"24\u00d722=528\", \"88\u00d713=1144\"],\n  [\"43\u00d748=2064\", \"67\u00d770=4690\"],\n  [\"60\u00d764=3840\", \"88\u00d784=7392\"],\n  [\"53\u00d758=3074\", \"38\u00d727=1026\"],\n  [\"65\u00d750=3250\", \"33\u00d745=1485\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    console.log(\"WARNING: not found -> \" + oldText);\n    continue;\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"87\u00d733=2871\", \"96\u00d721=2016\"),\n    @(\"60\u00d773=4380\", \"73\u00d788=6424\"),\n    @(\"51\u00d784=4284\", \"61\u00d742=2562\"),\n    @(\"84\u00d727=2268\", \"56\u00d779=4424\"),\n    @(\"67\u00d744=2948\", \"33\u00d712=396\"),\n    @(\"36\u00d751=1836\", \"68\u00d776=5168\"),\n    @(\"13\u00d737=481\", \"62\u00d745=2790\"),\n    @(\"41\u00d765=2665\", \"51\u00d739=1989\"),\n    @(\"81\u00d738=3078\", \"39\u00d780=3120\"),\n    @(\"51\u00d717=867\", \"53\u00d713=689\"),\n    @(\"31\u00d714=434\", \"72\u00d733=2376\"),\n    @(\"22\u00d744=968\", \"83\u00d759=4897\"),\n    @(\"43\u00d716=688\", \"49\u00d788=4312\"),\n    @(\"84\u00d784=7056\", \"66\u00d731=2046\"),\n    @(\"67\u00d775=5025\", \"16\u00d771=1136\"),\n    @(\"50\u00d713=650\", \"21\u00d774=1554\"),\n    @(\"87\u00d779=6873\", \"88\u00d764=5632\"),\n    @(\"33\u00d780=2640\", \"27\u00d789=2403\"),\n    @(\"91\u00d721=1911\", \"61\u00d723=1403\"),\n    @(\"64\u00d711=704\", \"84\u00d734=2856\"),\n    @(\"24\u00d722=528\", \"88\u00d713=1144\"),\n    @(\"43\u00d748=2064\", \"67\u00d770=4690\"),\n    @(\"60\u00d764=3840\", \"88\u00d784=7392\"),\n    @(\"53\u00d758=3074\", \"38\u00d727=1026\"),\n    @(\"65\u00d750=3250\", \"33\u00d745=1485\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        Write-Output \"WARNING: not found -> $oldText\"\n    }\n}"}
